$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# 1. Remove the "Engine" sheet (no longer needed)
$wb.Worksheets.Item("Engine").Delete() | Out-Null

# 2. "Air Volume" sheet: relabel rows - Compressor x -> VD x, "Compressor"/"Turned on" -> "VD"/"Volume"
$wsAir = $wb.Worksheets.Item("Air Volume")
$wsAir.Range("B1").Value = "Volume"
$wsAir.Range("A2").Value = "VD 1"
$wsAir.Range("A3").Value = "VD 2"
$wsAir.Range("A4").Value = "VD 3"
$wsAir.Range("A5").Value = "VD 4"
$wsAir.Range("A1").Value = "VD"

# 3. "Energy Consumption" sheet: append the next 7 time-series readings (rows 15-21)
$wsEnergy = $wb.Worksheets.Item("Energy Consumption")
$newRows = @(
  @(0.48263888888888901, 230, 380, 295, 0),
  @(0.48611111111111099, 200, 350, 295, 0),
  @(0.48958333333333298, 210, 320, 295, 0),
  @(0.49305555555555503, 200, 200, 295, 0),
  @(0.49652777777777701, 200, 150, 295, 0),
  @(0.499999999999999,   200, 100, 295, 0),
  @(0.50347222222222099, 220, 310, 295, 0)
)
$r = 15
foreach ($row in $newRows) {
  $wsEnergy.Range("A$r").Value = $row[0]
  $wsEnergy.Range("A$r").NumberFormat = "h:mm"
  $wsEnergy.Range("B$r").Value = $row[1]
  $wsEnergy.Range("C$r").Value = $row[2]
  $wsEnergy.Range("D$r").Value = $row[3]
  $wsEnergy.Range("E$r").Value = $row[4]
  $r++
}
$wsEnergy.Range("D24").Select() | Out-Null

# 4. Make "Air Volume" (now the first sheet) the active tab, default selection
$wsAir.Activate() | Out-Null
$wsAir.Range("A1").Select() | Out-Null
